# Auto-generated edit script applying the diff changes to sheet "Before FS-DR"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Before FS-DR")

# --- Simple value corrections (existing rows) ---
$ws.Range("B5").Value = 0.9804
$ws.Range("C5").Value = 0.7276
$ws.Range("B51").Value = 0.982
$ws.Range("C51").Value = 0.9811
$ws.Range("B63").Value = 0.9804
$ws.Range("C63").Value = 0.9792999999999999
$ws.Range("B84").Value = 0.982
$ws.Range("C84").Value = 0.7324000000000001
$ws.Range("B93").Value = 0.9757
$ws.Range("C93").Value = 0.722

# --- Newly filled rows (previously only had the Fold label in column A) ---
# Row 99
$ws.Range("B99").Value = 0.9748
$ws.Range("C99").Value = 0.8749
$ws.Range("D99").Value = 0.9721
$ws.Range("E99").Value = 0.9183
$ws.Range("F99").Value = 0.982
$ws.Range("G99").Value = 0.8662
$ws.Range("H99").Value = 0.9946
$ws.Range("I99").Value = 0.9767
$ws.Range("J99").Value = 0.9748
$ws.Range("K99").Value = 0.8946
$ws.Range("L99").Value = "C=1"
$ws.Range("M99").Value = "n_neighbors=5; weights=distance"
$ws.Range("N99").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O99").Value = "max_depth=None; n_estimators=200"
$ws.Range("P99").Value = "alpha=0.001; hidden_layer_sizes=(50, 50)"

# Row 100
$ws.Range("B100").Value = 0.9748
$ws.Range("C100").Value = 0.8282
$ws.Range("D100").Value = 0.973
$ws.Range("E100").Value = 0.8625
$ws.Range("F100").Value = 0.9856
$ws.Range("G100").Value = 0.8773
$ws.Range("H100").Value = 0.9901
$ws.Range("I100").Value = 0.9023
$ws.Range("J100").Value = 0.9775
$ws.Range("K100").Value = 0.8754999999999999
$ws.Range("L100").Value = "C=10"
$ws.Range("M100").Value = "n_neighbors=3; weights=distance"
$ws.Range("N100").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O100").Value = "max_depth=None; n_estimators=100"
$ws.Range("P100").Value = "alpha=0.0001; hidden_layer_sizes=(100,)"

# Row 101
$ws.Range("B101").Value = 0.9766
$ws.Range("C101").Value = 0.8999
$ws.Range("D101").Value = 0.9739
$ws.Range("E101").Value = 0.8841
$ws.Range("F101").Value = 0.9838
$ws.Range("G101").Value = 0.8391
$ws.Range("H101").Value = 0.9874000000000001
$ws.Range("I101").Value = 0.8847
$ws.Range("J101").Value = 0.982
$ws.Range("K101").Value = 0.9336
$ws.Range("L101").Value = "C=0.1"
$ws.Range("M101").Value = "n_neighbors=3; weights=distance"
$ws.Range("N101").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O101").Value = "max_depth=None; n_estimators=100"
$ws.Range("P101").Value = "alpha=0.0001; hidden_layer_sizes=(50,)"

# Row 102
$ws.Range("B102").Value = 0.964
$ws.Range("C102").Value = 0.849
$ws.Range("D102").Value = 0.9649
$ws.Range("E102").Value = 0.8491
$ws.Range("F102").Value = 0.9856
$ws.Range("G102").Value = 0.8944
$ws.Range("H102").Value = 0.9883
$ws.Range("I102").Value = 0.9320000000000001
$ws.Range("J102").Value = 0.973
$ws.Range("K102").Value = 0.8584000000000001
$ws.Range("L102").Value = "C=0.1"
$ws.Range("M102").Value = "n_neighbors=5; weights=distance"
$ws.Range("N102").Value = "max_depth=None; min_samples_split=5"
$ws.Range("O102").Value = "max_depth=None; n_estimators=200"
$ws.Range("P102").Value = "alpha=0.01; hidden_layer_sizes=(50, 50)"

# Row 103
$ws.Range("B103").Value = 0.9703000000000001
$ws.Range("C103").Value = 0.871
$ws.Range("D103").Value = 0.9685
$ws.Range("E103").Value = 0.8649
$ws.Range("F103").Value = 0.9811
$ws.Range("G103").Value = 0.9031
$ws.Range("H103").Value = 0.9955000000000001
$ws.Range("I103").Value = 0.9717
$ws.Range("J103").Value = 0.9802
$ws.Range("K103").Value = 0.9115
$ws.Range("L103").Value = "C=0.1"
$ws.Range("M103").Value = "n_neighbors=3; weights=distance"
$ws.Range("N103").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O103").Value = "max_depth=None; n_estimators=200"
$ws.Range("P103").Value = "alpha=0.001; hidden_layer_sizes=(100,)"

# Row 104
$ws.Range("B104").Value = 0.9721
$ws.Range("C104").Value = 0.8941
$ws.Range("D104").Value = 0.9792999999999999
$ws.Range("E104").Value = 0.8935
$ws.Range("F104").Value = 0.982
$ws.Range("G104").Value = 0.8677
$ws.Range("H104").Value = 0.9892
$ws.Range("I104").Value = 0.9022
$ws.Range("J104").Value = 0.9865
$ws.Range("K104").Value = 0.9394
$ws.Range("L104").Value = "C=0.1"
$ws.Range("M104").Value = "n_neighbors=5; weights=distance"
$ws.Range("N104").Value = "max_depth=20; min_samples_split=2"
$ws.Range("O104").Value = "max_depth=None; n_estimators=200"
$ws.Range("P104").Value = "alpha=0.01; hidden_layer_sizes=(100,)"

# Row 105
$ws.Range("B105").Value = 0.9622000000000001
$ws.Range("C105").Value = 0.8342000000000001
$ws.Range("D105").Value = 0.9667
$ws.Range("E105").Value = 0.8762
$ws.Range("F105").Value = 0.9802
$ws.Range("G105").Value = 0.8566
$ws.Range("H105").Value = 0.9901
$ws.Range("I105").Value = 0.9261
$ws.Range("J105").Value = 0.9792999999999999
$ws.Range("K105").Value = 0.8972
$ws.Range("L105").Value = "C=0.1"
$ws.Range("M105").Value = "n_neighbors=3; weights=distance"
$ws.Range("N105").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O105").Value = "max_depth=20; n_estimators=100"
$ws.Range("P105").Value = "alpha=0.001; hidden_layer_sizes=(50, 50)"

# Row 106
$ws.Range("B106").Value = 0.9622000000000001
$ws.Range("C106").Value = 0.7895
$ws.Range("D106").Value = 0.9694
$ws.Range("E106").Value = 0.802
$ws.Range("F106").Value = 0.9838
$ws.Range("G106").Value = 0.8246
$ws.Range("H106").Value = 0.991
$ws.Range("I106").Value = 0.8822
$ws.Range("J106").Value = 0.9694
$ws.Range("K106").Value = 0.7892
$ws.Range("L106").Value = "C=1"
$ws.Range("M106").Value = "n_neighbors=3; weights=distance"
$ws.Range("N106").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O106").Value = "max_depth=20; n_estimators=200"
$ws.Range("P106").Value = "alpha=0.001; hidden_layer_sizes=(100,)"

# Row 107
$ws.Range("B107").Value = 0.9676
$ws.Range("C107").Value = 0.8258
$ws.Range("D107").Value = 0.9631
$ws.Range("E107").Value = 0.8478
$ws.Range("F107").Value = 0.9847
$ws.Range("G107").Value = 0.8679
$ws.Range("H107").Value = 0.9937
$ws.Range("I107").Value = 0.9515
$ws.Range("J107").Value = 0.9685
$ws.Range("K107").Value = 0.8401999999999999
$ws.Range("L107").Value = "C=0.1"
$ws.Range("M107").Value = "n_neighbors=3; weights=distance"
$ws.Range("N107").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O107").Value = "max_depth=20; n_estimators=200"
$ws.Range("P107").Value = "alpha=0.01; hidden_layer_sizes=(100,)"

# Row 108
$ws.Range("B108").Value = 0.9757
$ws.Range("C108").Value = 0.8705000000000001
$ws.Range("D108").Value = 0.9685
$ws.Range("E108").Value = 0.8767
$ws.Range("F108").Value = 0.9865
$ws.Range("G108").Value = 0.9156
$ws.Range("H108").Value = 0.9892
$ws.Range("I108").Value = 0.9316
$ws.Range("J108").Value = 0.9784
$ws.Range("K108").Value = 0.9147999999999999
$ws.Range("L108").Value = "C=1"
$ws.Range("M108").Value = "n_neighbors=9; weights=distance"
$ws.Range("N108").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O108").Value = "max_depth=None; n_estimators=200"
$ws.Range("P108").Value = "alpha=0.001; hidden_layer_sizes=(50, 50)"

# Row 111
$ws.Range("B111").Value = 0.9397
$ws.Range("C111").Value = 0.5451
$ws.Range("D111").Value = 0.9406
$ws.Range("E111").Value = 0.5575
$ws.Range("F111").Value = 0.9478
$ws.Range("G111").Value = 0.532
$ws.Range("H111").Value = 0.9568
$ws.Range("I111").Value = 0.5848
$ws.Range("J111").Value = 0.9451000000000001
$ws.Range("K111").Value = 0.5782
$ws.Range("L111").Value = "C=0.1"
$ws.Range("M111").Value = "n_neighbors=9; weights=distance"
$ws.Range("N111").Value = "max_depth=20; min_samples_split=2"
$ws.Range("O111").Value = "max_depth=None; n_estimators=100"
$ws.Range("P111").Value = "alpha=0.0001; hidden_layer_sizes=(50,)"

# Row 112
$ws.Range("B112").Value = 0.9352
$ws.Range("C112").Value = 0.5292
$ws.Range("D112").Value = 0.9469
$ws.Range("E112").Value = 0.523
$ws.Range("F112").Value = 0.9451000000000001
$ws.Range("G112").Value = 0.4759
$ws.Range("H112").Value = 0.9559
$ws.Range("I112").Value = 0.5591
$ws.Range("J112").Value = 0.9343
$ws.Range("K112").Value = 0.499
$ws.Range("L112").Value = "C=0.1"
$ws.Range("M112").Value = "n_neighbors=5; weights=distance"
$ws.Range("N112").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O112").Value = "max_depth=None; n_estimators=100"
$ws.Range("P112").Value = "alpha=0.001; hidden_layer_sizes=(50,)"

# Row 113
$ws.Range("B113").Value = 0.9370000000000001
$ws.Range("C113").Value = 0.5281
$ws.Range("D113").Value = 0.9397
$ws.Range("E113").Value = 0.5224
$ws.Range("F113").Value = 0.9469
$ws.Range("G113").Value = 0.4974
$ws.Range("H113").Value = 0.9523
$ws.Range("I113").Value = 0.542
$ws.Range("J113").Value = 0.946
$ws.Range("K113").Value = 0.5621
$ws.Range("L113").Value = "C=0.1"
$ws.Range("M113").Value = "n_neighbors=3; weights=distance"
$ws.Range("N113").Value = "max_depth=20; min_samples_split=10"
$ws.Range("O113").Value = "max_depth=None; n_estimators=200"
$ws.Range("P113").Value = "alpha=0.0001; hidden_layer_sizes=(50,)"

# Row 114
$ws.Range("B114").Value = 0.9361
$ws.Range("C114").Value = 0.4686
$ws.Range("D114").Value = 0.9433
$ws.Range("E114").Value = 0.4925
$ws.Range("F114").Value = 0.9496
$ws.Range("G114").Value = 0.5336
$ws.Range("H114").Value = 0.9496
$ws.Range("I114").Value = 0.4824
$ws.Range("J114").Value = 0.9433
$ws.Range("K114").Value = 0.5067
$ws.Range("L114").Value = "C=0.1"
$ws.Range("M114").Value = "n_neighbors=7; weights=distance"
$ws.Range("N114").Value = "max_depth=None; min_samples_split=10"
$ws.Range("O114").Value = "max_depth=10; n_estimators=100"
$ws.Range("P114").Value = "alpha=0.001; hidden_layer_sizes=(100,)"

# Row 115
$ws.Range("B115").Value = 0.9442
$ws.Range("C115").Value = 0.5275
$ws.Range("D115").Value = 0.9442
$ws.Range("E115").Value = 0.5492
$ws.Range("F115").Value = 0.9397
$ws.Range("G115").Value = 0.4806
$ws.Range("H115").Value = 0.9568
$ws.Range("I115").Value = 0.5658
$ws.Range("J115").Value = 0.9496
$ws.Range("K115").Value = 0.5319
$ws.Range("L115").Value = "C=0.1"
$ws.Range("M115").Value = "n_neighbors=7; weights=distance"
$ws.Range("N115").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O115").Value = "max_depth=None; n_estimators=100"
$ws.Range("P115").Value = "alpha=0.0001; hidden_layer_sizes=(100,)"

# Row 116
$ws.Range("B116").Value = 0.9352
$ws.Range("C116").Value = 0.5248
$ws.Range("D116").Value = 0.9379
$ws.Range("E116").Value = 0.5406
$ws.Range("F116").Value = 0.9478
$ws.Range("G116").Value = 0.5667
$ws.Range("H116").Value = 0.9514
$ws.Range("I116").Value = 0.5499000000000001
$ws.Range("J116").Value = 0.9370000000000001
$ws.Range("K116").Value = 0.552
$ws.Range("L116").Value = "C=0.01"
$ws.Range("M116").Value = "n_neighbors=5; weights=distance"
$ws.Range("N116").Value = "max_depth=20; min_samples_split=10"
$ws.Range("O116").Value = "max_depth=20; n_estimators=200"
$ws.Range("P116").Value = "alpha=0.001; hidden_layer_sizes=(100,)"

# Row 117
$ws.Range("B117").Value = 0.9559
$ws.Range("C117").Value = 0.6776
$ws.Range("D117").Value = 0.9505
$ws.Range("E117").Value = 0.5947
$ws.Range("F117").Value = 0.9523
$ws.Range("G117").Value = 0.6012999999999999
$ws.Range("H117").Value = 0.9595
$ws.Range("I117").Value = 0.6076
$ws.Range("J117").Value = 0.9577
$ws.Range("K117").Value = 0.6384
$ws.Range("L117").Value = "C=0.1"
$ws.Range("M117").Value = "n_neighbors=3; weights=distance"
$ws.Range("N117").Value = "max_depth=20; min_samples_split=10"
$ws.Range("O117").Value = "max_depth=20; n_estimators=200"
$ws.Range("P117").Value = "alpha=0.001; hidden_layer_sizes=(100,)"

# Row 118
$ws.Range("B118").Value = 0.9324
$ws.Range("C118").Value = 0.5308
$ws.Range("D118").Value = 0.9252
$ws.Range("E118").Value = 0.4419
$ws.Range("F118").Value = 0.9387
$ws.Range("G118").Value = 0.4699
$ws.Range("H118").Value = 0.9468
$ws.Range("I118").Value = 0.5046
$ws.Range("J118").Value = 0.9351
$ws.Range("K118").Value = 0.5127
$ws.Range("L118").Value = "C=0.1"
$ws.Range("M118").Value = "n_neighbors=7; weights=distance"
$ws.Range("N118").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O118").Value = "max_depth=20; n_estimators=100"
$ws.Range("P118").Value = "alpha=0.0001; hidden_layer_sizes=(100,)"

# Row 119
$ws.Range("B119").Value = 0.9351
$ws.Range("C119").Value = 0.5647
$ws.Range("D119").Value = 0.9423
$ws.Range("E119").Value = 0.5291
$ws.Range("F119").Value = 0.9514
$ws.Range("G119").Value = 0.5363
$ws.Range("H119").Value = 0.9523
$ws.Range("I119").Value = 0.5322
$ws.Range("J119").Value = 0.9459
$ws.Range("K119").Value = 0.5639999999999999
$ws.Range("L119").Value = "C=0.1"
$ws.Range("M119").Value = "n_neighbors=5; weights=distance"
$ws.Range("N119").Value = "max_depth=20; min_samples_split=10"
$ws.Range("O119").Value = "max_depth=20; n_estimators=200"
$ws.Range("P119").Value = "alpha=0.0001; hidden_layer_sizes=(50, 50)"

# Row 120
$ws.Range("B120").Value = 0.9342
$ws.Range("C120").Value = 0.5387999999999999
$ws.Range("D120").Value = 0.9387
$ws.Range("E120").Value = 0.514
$ws.Range("F120").Value = 0.9405
$ws.Range("G120").Value = 0.4756
$ws.Range("H120").Value = 0.9486
$ws.Range("I120").Value = 0.5158
$ws.Range("J120").Value = 0.9324
$ws.Range("K120").Value = 0.4974
$ws.Range("L120").Value = "C=0.1"
$ws.Range("M120").Value = "n_neighbors=5; weights=distance"
$ws.Range("N120").Value = "max_depth=None; min_samples_split=2"
$ws.Range("O120").Value = "max_depth=None; n_estimators=200"
$ws.Range("P120").Value = "alpha=0.001; hidden_layer_sizes=(50,)"

